$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix trailing whitespace in existing beach name strings
$ws.Range("B18").Value = "Port Glasgow"
$ws.Range("B36").Value = "NEW BEACH THAT I CREATED"

# Add new row 37 for the "BAWS" beach entry
$ws.Range("A37").Value = 105
$ws.Range("B37").Value = "BAWS"
$ws.Range("C37").Value = "BAWS"
$ws.Range("D37").Value = "BAWS"
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 3

$ws.Range("F37").Select()
